# "Rename 'wire_transfers' sheet to 'currency conversion to EUR'"
#
# The old "Wire transfer" label was replaced everywhere it appears (it is
# the text behind a shared string used by two rows in the "Fees" sheet),
# a couple of exchange-rate figures on "Foreign Currencies" were corrected,
# the "Fees" sheet's label column was widened to fit the new, longer text,
# and the ELSTER summary's "Fremdwaehrungen" total was refreshed to match.

$wb = $excel.ActiveWorkbook

# --- Fees: "Wire transfer" -> "Currency conversion or wire transfer" ------
# Both rows that used to read "Wire transfer" (B6 and B11) share the same
# underlying text, so both need to be updated for the rename to take full
# effect everywhere the old label showed up.
$feesWs = $wb.Worksheets.Item("Fees")
$newLabel = "Currency conversion or wire transfer"
$feesWs.Range("B6").Value = $newLabel
$feesWs.Range("B11").Value = $newLabel

# The label column (B) now holds much longer text, so widen it to fit --
# mirrors what Excel does to a best-fit column once its content grows.
$feesWs.Columns.Item(2).ColumnWidth = 32.33

# --- Foreign Currencies: corrected GOOG buy amount / EUR gain -------------
$fxWs = $wb.Worksheets.Item("Foreign Currencies")
$fxWs.Range("B7").Value = 155
$fxWs.Range("G7").Value = -10.64

# --- ELSTER - Summary: updated "Gewinn/Verlust aus Fremdwaehrungen" -------
$elsterWs = $wb.Worksheets.Item("ELSTER - Summary")
$elsterWs.Range("C7").Value = 67.67
